# Fruta / hortaliza, semanal
#
# The data rows (2-23) of the "Haba" sheet got their weekly price/volume
# tuples (Fecha + Volumen/Precio minimo/Precio maximo/Precio promedio
# ponderado/Precio $/Kg, i.e. columns D,J,K,L,M,P) reshuffled across rows -
# each row now carries the tuple that used to belong to a different row
# (row 4 keeps its own tuple). Everything else (Mercado, Region, Categoria,
# Variedad, Calidad, Unidad, Origen, Clasificacion, etc.) stays put.
#
# Map: destination row -> source row (where its "new" tuple used to live).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2  = 23
    3  = 5
    4  = 4
    5  = 9
    6  = 10
    7  = 20
    8  = 14
    9  = 3
    10 = 2
    11 = 16
    12 = 13
    13 = 11
    14 = 7
    15 = 12
    16 = 15
    17 = 19
    18 = 21
    19 = 6
    20 = 8
    21 = 22
    22 = 17
    23 = 18
}

# Columns that travel together as one tuple: Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg.
$cols = @(4, 10, 11, 12, 13, 16)

# Snapshot every source row's tuple BEFORE writing anything, since this is a
# full permutation (every row is simultaneously a source and a destination).
$snapshot = @{}
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    if (-not $snapshot.ContainsKey($srcRow)) {
        $vals = @{}
        foreach ($c in $cols) {
            $vals[$c] = $ws.Cells.Item($srcRow, $c).Value2
        }
        $snapshot[$srcRow] = $vals
    }
}

# Now apply the snapshot to each destination row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    if ($srcRow -eq $destRow) {
        continue
    }
    $vals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $vals[$c]
    }
}
